$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    ,@("2023-03-22-14:00", 28269.7, 28758.1, 28264.6, 28546, 28594.01502305624, 28594.01502305624, 97107.311)
    ,@("2023-03-22-15:00", 28545.9, 28731, 28423, 28586.2, 28594.01502305624, 28670.33106257414, 65870.58100000001)
    ,@("2023-03-22-16:00", 28586.1, 28696.3, 28469.1, 28636.6, 28670.33106257414, 28722.41414103485, 30881.951)
    ,@("2023-03-22-17:00", 28636.7, 28639.1, 28283.5, 28500, 28722.41414103485, 28646.12986207716, 57523.283)
    ,@("2023-03-22-18:00", 28500.4, 28881, 27769, 28017.8, 28646.12986207716, 28332.9595793561, 208685.849)
    ,@("2023-03-22-19:00", 28018.9, 28020.3, 26620, 26654.1, 28332.9595793561, 27093.8079718299, 215182.357)
    ,@("2023-03-22-20:00", 26653.3, 27372, 26590, 27295, 27093.8079718299, 27421.03296262654, 82764.427)
    ,@("2023-03-22-21:00", 27294.9, 27425.3, 27056.4, 27065.3, 27421.03296262654, 26951.71434903986, 40658.496)
    ,@("2023-03-22-22:00", 27065.3, 27224, 26956.3, 27189.2, 26951.71434903986, 27108.25674876536, 26901.703)
    ,@("2023-03-22-23:00", 27189.2, 27333, 27169.1, 27233.8, 27108.25674876536, 27247.31263966346, 18277.151)
    ,@("2023-03-23-00:00", 27233.8, 27421.8, 27178, 27305.6, 27247.31263966346, 27161.47371213962, 23014.289)
    ,@("2023-03-23-01:00", 27305.5, 27309.2, 27101, 27137.9, 27161.47371213962, 27101.12690264694, 15997.838)
    ,@("2023-03-23-02:00", 27138, 27338.3, 27130.5, 27336.7, 27101.12690264694, 27214.7209693265, 11407.743)
    ,@("2023-03-23-03:00", 27336.6, 27415, 27282.1, 27343.7, 27214.7209693265, 27577.73025693072, 13425.824)
    ,@("2023-03-23-04:00", 27343.7, 27411, 27231.3, 27262.4, 27577.73025693072, 27423.21722348171, 13170.433)
    ,@("2023-03-23-05:00", 27262.5, 27680, 27262.1, 27608, 27423.21722348171, 27571.58559090094, 35976.002)
    ,@("2023-03-23-06:00", 27608, 27675.4, 27540, 27632.9, 27571.58559090094, 27587.91775878609, 20009.877)
    ,@("2023-03-23-07:00", 27633, 27660.6, 27518.8, 27616.4, 27587.91775878609, 27591.39579847406, 15381.825)
)

$startRow = 63
for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $rowData = $data[$i]
    $ws.Cells.Item($r, 1).Value = $rowData[0]
    for ($col = 2; $col -le 8; $col++) {
        $ws.Cells.Item($r, $col).Value = [double]$rowData[$col - 1]
    }
}

$wb.Save()
